# Append the new staff numbers (rows 6-11) below the existing data on Sheet1,
# then move the active selection to the new last cell (A11), matching the
# "after" state of xl/worksheets/sheet1.xml in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(222005, 222006, 222101, 222103, 222104, 222303)

$row = 6
foreach ($v in $newValues) {
    $ws.Cells.Item($row, 1).Value = $v
    $row++
}

# Move selection to the new last populated cell (A11), mirroring the
# workbook's saved selection/activeCell state after the edit.
$ws.Range("A11").Select() | Out-Null
